$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logins")

$values = @("provider11","provider12","provider13","provider14","provider15","provider16","provider17","provider18","provider19","provider20","provider21")

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
